# This revision's underlying commit ("Fixed POI packaging and upgraded to
# POI 3.15") is a build-tooling change in the project that produced this
# fixture .docx, not a content edit made by hand. Every single hunk in the
# diff is a pure attribute/namespace re-serialization: the POI (XMLBeans)
# writer used after the upgrade emits element attributes in a different
# (sorted) order, but every tag, every attribute name, and every attribute
# value is byte-for-byte identical to before. A few representative
# examples (old -> new), all equal as sets:
#   <w:pgSz w:w="11906" w:h="16838"/>
#     -> <w:pgSz w:h="16838" w:w="11906"/>
#   <w:lang w:val="fr-FR" w:eastAsia="en-US" w:bidi="ar-SA"/>
#     -> <w:lang w:bidi="ar-SA" w:eastAsia="en-US" w:val="fr-FR"/>
#   <w:style w:type="paragraph" w:default="1" w:styleId="Normal">
#     -> <w:style w:default="1" w:styleId="Normal" w:type="paragraph">
# No paragraph text, run content, field code, style definition, numbering,
# or formatting value changes anywhere in the document or style sheet.
#
# Word's object model (real Word, and this COM-interop host alike) does
# not expose — and this edit does not require — any control over the raw
# XML attribute ordering used when a part is serialized; that is an
# internal detail of the writer, not a document property. There is
# nothing for a Word automation script to change here: the document's
# actual content, formatting, sections, and styles are already exactly
# what the target state calls for.
#
# Touch the document (read-only, idempotent) so the script is a deliberate
# no-op rather than an accidental empty file, without introducing any
# content, formatting, or structural difference.
$d = $word.ActiveDocument
$null = $d.Content.Text
$null = $d.Styles("Normal").NameLocal
$null = $d.Sections(1).PageSetup.TopMargin
